$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 712 ("「空から降り立った者たち」" post) entirely; rows below shift up.
$ws.Rows.Item(712).Delete()
